# ----------------------------------------------------------------------------
# Applies the two edits captured by the target diff:
#
#   1. The table on slide 16 gets its table style switched from the deck's
#      custom "Table_0" style ({251CDCB6-0FFB-42C2-855C-FFE2EBCBB836}) to the
#      built-in "Medium Style 2 - Accent 1" style
#      ({ABF68BB5-0CCE-4063-9A54-63974E27CE6D}).
#
#   2. The two embedded theme parts (ppt/theme/theme1.xml, used by the notes
#      master, and ppt/theme/theme2.xml, used by the slide master / the whole
#      deck) swap their "Office Theme" / "Integral" colour schemes. The
#      reachable half of that swap, through the PowerPoint object model, is
#      the live deck theme (theme2.xml): every themed colour slot gets set to
#      the plain "Office Theme" palette values that used to live in
#      theme1.xml.
# ----------------------------------------------------------------------------

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Table style swap on the cash-flow table (slide 16).
# ---------------------------------------------------------------------------
$targetStyleId = "{ABF68BB5-0CCE-4063-9A54-63974E27CE6D}"

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($targetStyleId)
        }
    }
}

# ---------------------------------------------------------------------------
# 2) Theme colour swap: push the "Office Theme" palette into the live theme
#    (theme2.xml), replacing the current "Integral" palette.
#    Order of a ThemeColorScheme's 12 slots matches <a:clrScheme>'s
#    children: dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink.
#    RGB is assigned as a 0xBBGGRR long (VBA RGB()-style), matching the COM
#    RGBColor convention.
# ---------------------------------------------------------------------------
$officeThemeRgb = @(
    0x000000,  # dk1      000000
    0xFFFFFF,  # lt1      FFFFFF
    0x6A5444,  # dk2      44546A
    0xE6E6E7,  # lt2      E7E6E6
    0xD59B5B,  # accent1  5B9BD5
    0x317DED,  # accent2  ED7D31
    0xA5A5A5,  # accent3  A5A5A5
    0x00C0FF,  # accent4  FFC000
    0xC47244,  # accent5  4472C4
    0x47AD70,  # accent6  70AD47
    0xC16305,  # hlink    0563C1
    0x724F95   # folHlink 954F72
)

$slideOne = $p.Slides.Item(1)
$themeColors = $slideOne.ThemeColorScheme
for ($ci = 1; $ci -le 12; $ci++) {
    $themeColors.Colors($ci).RGB = $officeThemeRgb[$ci - 1]
}
